$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean sheet so the used range / dimension exactly matches the target
$ws.Cells.Clear()

# --- Date style helper: format used by the header row (row 1) and the anchor column (A) ---
function Set-DateStyle($rng) {
    $rng.NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $rng.Font.Bold = $true
    $rng.HorizontalAlignment = -4108
    $rng.VerticalAlignment = -4160
    $rng.Borders.LineStyle = 1
}

# --- Row 1: period-end dates across columns B:AZ ---
$row1Dates = @(39583,39765,39948,40130,40310,40494,40676,40862,41044,41228,41409,41592,41774,41957,42137,42321,42503,42689,42867,43053,43145,43235,43326,43418,43510,43600,43691,43783,43875,43966,44068,44159,44251,44341,44432,44525,44617,44706,44798,44890,44981,45071,45163,45254,45345,45436,45534,45618,45713,45800,45891)
for ($i = 0; $i -lt $row1Dates.Length; $i++) {
    $cell = $ws.Cells.Item(1, $i + 2)
    Set-DateStyle $cell
    $cell.Value = $row1Dates[$i]
}

# --- Column A: anchor (vintage) dates across rows 2:22 ---
$colADates = @(39813,40178,40543,40908,41274,41639,42004,42369,42735,43100,43465,43830,44196,44561,44926,45291,45657,46022,46387,46752,47118)
for ($i = 0; $i -lt $colADates.Length; $i++) {
    $cell = $ws.Cells.Item($i + 2, 1)
    Set-DateStyle $cell
    $cell.Value = $colADates[$i]
}

# --- Data matrix: one contiguous run of forecast values per row, columns grow with the row ---
$r3 = @(5.515135932605308,2.741420124134053,-2.764179109379705,-3.872359107260159,-3.872359107260159,-3.872359107260159,-3.872359107260159,-3.872359107260159,-3.872359107260159,-3.872359107260159,-3.872359107260159,-3.872359107260159,-3.872359107260159,-3.872359107260159,-3.872359107260159,-3.872359107260159,-3.872359107260159,-3.872359107260159,-3.872359107260159,-3.872359107260159,-3.872359107260159,-3.872359107260159,-3.872359107260159,-3.872359107260159,-3.872359107260159,-3.872359107260159,-3.872359107260159,-3.872359107260159,-3.872359107260159,-3.872359107260159,-3.872359107260159,-3.872359107260159,-3.872359107260159,-3.872359107260159,-3.872359107260159,-3.872359107260159,-3.872359107260159,-3.872359107260159,-3.872359107260159,-3.872359107260159,-3.872359107260159,-3.872359107260159,-3.872359107260159,-3.872359107260159,-3.872359107260159,-3.872359107260159,-3.872359107260159,-3.872359107260159,-3.872359107260159,-3.872359107260159,-3.872359107260159)
for ($i = 0; $i -lt $r3.Length; $i++) { $ws.Cells.Item(3, 2 + $i).Value = $r3[$i] }
$r4 = @(5.566598885960738,3.4712032095223,-0.2596132895792413,-1.796159642965267,3.38265053313096,4.530477057343663,4.530477057343663,4.530477057343663,4.530477057343663,4.530477057343663,4.530477057343663,4.530477057343663,4.530477057343663,4.530477057343663,4.530477057343663,4.530477057343663,4.530477057343663,4.530477057343663,4.530477057343663,4.530477057343663,4.530477057343663,4.530477057343663,4.530477057343663,4.530477057343663,4.530477057343663,4.530477057343663,4.530477057343663,4.530477057343663,4.530477057343663,4.530477057343663,4.530477057343663,4.530477057343663,4.530477057343663,4.530477057343663,4.530477057343663,4.530477057343663,4.530477057343663,4.530477057343663,4.530477057343663,4.530477057343663,4.530477057343663,4.530477057343663,4.530477057343663,4.530477057343663,4.530477057343663,4.530477057343663,4.530477057343663,4.530477057343663,4.530477057343663,4.530477057343663,4.530477057343663)
for ($i = 0; $i -lt $r4.Length; $i++) { $ws.Cells.Item(4, 2 + $i).Value = $r4[$i] }
$r5 = @(-0.06191106053077267,-1.072567702512639,1.421687849828523,3.372216600811506,5.813045170083808,6.833902841285977,6.833902841285977,6.833902841285977,6.833902841285977,6.833902841285977,6.833902841285977,6.833902841285977,6.833902841285977,6.833902841285977,6.833902841285977,6.833902841285977,6.833902841285977,6.833902841285977,6.833902841285977,6.833902841285977,6.833902841285977,6.833902841285977,6.833902841285977,6.833902841285977,6.833902841285977,6.833902841285977,6.833902841285977,6.833902841285977,6.833902841285977,6.833902841285977,6.833902841285977,6.833902841285977,6.833902841285977,6.833902841285977,6.833902841285977,6.833902841285977,6.833902841285977,6.833902841285977,6.833902841285977,6.833902841285977,6.833902841285977,6.833902841285977,6.833902841285977,6.833902841285977,6.833902841285977,6.833902841285977,6.833902841285977,6.833902841285977,6.833902841285977)
for ($i = 0; $i -lt $r5.Length; $i++) { $ws.Cells.Item(5, 4 + $i).Value = $r5[$i] }
$r6 = @(1.047976847023846,1.316610894175119,2.371160938652705,5.794553957309168,4.597183386292891,4.166536506645224,4.166536506645224,4.166536506645224,4.166536506645224,4.166536506645224,4.166536506645224,4.166536506645224,4.166536506645224,4.166536506645224,4.166536506645224,4.166536506645224,4.166536506645224,4.166536506645224,4.166536506645224,4.166536506645224,4.166536506645224,4.166536506645224,4.166536506645224,4.166536506645224,4.166536506645224,4.166536506645224,4.166536506645224,4.166536506645224,4.166536506645224,4.166536506645224,4.166536506645224,4.166536506645224,4.166536506645224,4.166536506645224,4.166536506645224,4.166536506645224,4.166536506645224,4.166536506645224,4.166536506645224,4.166536506645224,4.166536506645224,4.166536506645224,4.166536506645224,4.166536506645224,4.166536506645224,4.166536506645224,4.166536506645224)
for ($i = 0; $i -lt $r6.Length; $i++) { $ws.Cells.Item(6, 6 + $i).Value = $r6[$i] }
$r7 = @(2.253457754589228,5.18838101923762,6.213583554874536,4.453225076234824,2.113053977048707,2.669880057548091,2.669880057548091,2.669880057548091,2.669880057548091,2.669880057548091,2.669880057548091,2.669880057548091,2.669880057548091,2.669880057548091,2.669880057548091,2.669880057548091,2.669880057548091,2.669880057548091,2.669880057548091,2.669880057548091,2.669880057548091,2.669880057548091,2.669880057548091,2.669880057548091,2.669880057548091,2.669880057548091,2.669880057548091,2.669880057548091,2.669880057548091,2.669880057548091,2.669880057548091,2.669880057548091,2.669880057548091,2.669880057548091,2.669880057548091,2.669880057548091,2.669880057548091,2.669880057548091,2.669880057548091,2.669880057548091,2.669880057548091,2.669880057548091,2.669880057548091,2.669880057548091,2.669880057548091)
for ($i = 0; $i -lt $r7.Length; $i++) { $ws.Cells.Item(7, 8 + $i).Value = $r7[$i] }
$r8 = @(6.4658713724427,5.226878449528161,3.80132608443593,4.443665867214519,5.262295419893648,5.50293301232252,5.50293301232252,5.50293301232252,5.50293301232252,5.50293301232252,5.50293301232252,5.50293301232252,5.50293301232252,5.50293301232252,5.50293301232252,5.50293301232252,5.50293301232252,5.50293301232252,5.50293301232252,5.50293301232252,5.50293301232252,5.50293301232252,5.50293301232252,5.50293301232252,5.50293301232252,5.50293301232252,5.50293301232252,5.50293301232252,5.50293301232252,5.50293301232252,5.50293301232252,5.50293301232252,5.50293301232252,5.50293301232252,5.50293301232252,5.50293301232252,5.50293301232252,5.50293301232252,5.50293301232252,5.50293301232252,5.50293301232252,5.50293301232252,5.50293301232252)
for ($i = 0; $i -lt $r8.Length; $i++) { $ws.Cells.Item(8, 10 + $i).Value = $r8[$i] }
$r9 = @(3.947168334316919,3.983760133672121,3.911207409579687,4.93942156312468,4.78031692483154,4.829481320500406,4.829481320500406,4.829481320500406,4.829481320500406,4.829481320500406,4.829481320500406,4.829481320500406,4.829481320500406,4.829481320500406,4.829481320500406,4.829481320500406,4.829481320500406,4.829481320500406,4.829481320500406,4.829481320500406,4.829481320500406,4.829481320500406,4.829481320500406,4.829481320500406,4.829481320500406,4.829481320500406,4.829481320500406,4.829481320500406,4.829481320500406,4.829481320500406,4.829481320500406,4.829481320500406,4.829481320500406,4.829481320500406,4.829481320500406,4.829481320500406,4.829481320500406,4.829481320500406,4.829481320500406,4.829481320500406,4.829481320500406)
for ($i = 0; $i -lt $r9.Length; $i++) { $ws.Cells.Item(9, 12 + $i).Value = $r9[$i] }
$r10 = @(3.718405228598609,4.397175883931581,4.481753591536197,5.304093839766777,5.386945580119185,5.100281927437122,5.100281927437122,5.100281927437122,5.100281927437122,5.100281927437122,5.100281927437122,5.100281927437122,5.100281927437122,5.100281927437122,5.100281927437122,5.100281927437122,5.100281927437122,5.100281927437122,5.100281927437122,5.100281927437122,5.100281927437122,5.100281927437122,5.100281927437122,5.100281927437122,5.100281927437122,5.100281927437122,5.100281927437122,5.100281927437122,5.100281927437122,5.100281927437122,5.100281927437122,5.100281927437122,5.100281927437122,5.100281927437122,5.100281927437122,5.100281927437122,5.100281927437122,5.100281927437122,5.100281927437122)
for ($i = 0; $i -lt $r10.Length; $i++) { $ws.Cells.Item(10, 14 + $i).Value = $r10[$i] }
$r11 = @(4.481643758166953,5.391718764442932,5.291601650505706,4.773858406282372,4.884223728030879,5.161358932333737,5.161358932333737,5.161358932333737,5.161358932333737,5.161358932333737,5.161358932333737,5.161358932333737,5.161358932333737,5.161358932333737,5.161358932333737,5.161358932333737,5.161358932333737,5.161358932333737,5.161358932333737,5.161358932333737,5.161358932333737,5.161358932333737,5.161358932333737,5.161358932333737,5.161358932333737,5.161358932333737,5.161358932333737,5.161358932333737,5.161358932333737,5.161358932333737,5.161358932333737,5.161358932333737,5.161358932333737,5.161358932333737,5.161358932333737,5.161358932333737,5.161358932333737)
for ($i = 0; $i -lt $r11.Length; $i++) { $ws.Cells.Item(11, 16 + $i).Value = $r11[$i] }
$r12 = @(5.256062028567232,4.867860541745572,4.848884192354119,5.448823958443616,6.111775080355519,6.072005530313129,5.969015988310433,5.902681694119694,5.902681694119694,5.902681694119694,5.902681694119694,5.902681694119694,5.902681694119694,5.902681694119694,5.902681694119694,5.902681694119694,5.902681694119694,5.902681694119694,5.902681694119694,5.902681694119694,5.902681694119694,5.902681694119694,5.902681694119694,5.902681694119694,5.902681694119694,5.902681694119694,5.902681694119694,5.902681694119694,5.902681694119694,5.902681694119694,5.902681694119694,5.902681694119694,5.902681694119694,5.902681694119694,5.902681694119694)
for ($i = 0; $i -lt $r12.Length; $i++) { $ws.Cells.Item(12, 18 + $i).Value = $r12[$i] }
$r13 = @(4.819877891111335,5.22223210260726,5.5540247491221,5.47980442665561,5.321568540215016,5.03659417252571,3.94180168325462,4.927928448556984,3.961704997901161,3.884502719230132,3.884502719230132,3.884502719230132,3.884502719230132,3.884502719230132,3.884502719230132,3.884502719230132,3.884502719230132,3.884502719230132,3.884502719230132,3.884502719230132,3.884502719230132,3.884502719230132,3.884502719230132,3.884502719230132,3.884502719230132,3.884502719230132,3.884502719230132,3.884502719230132,3.884502719230132,3.884502719230132,3.884502719230132,3.884502719230132,3.884502719230132)
for ($i = 0; $i -lt $r13.Length; $i++) { $ws.Cells.Item(13, 20 + $i).Value = $r13[$i] }
$r14 = @(5.484025951319293,5.40641359612255,5.259836672320484,5.077291536394846,5.802765867180804,4.342221082693465,4.009670676786059,3.818894565497888,2.502519143054571,-3.840397826549158,-3.840397826549158,-3.840397826549158,-3.840397826549158,-3.840397826549158,-3.840397826549158,-3.840397826549158,-3.840397826549158,-3.840397826549158,-3.840397826549158,-3.840397826549158,-3.840397826549158,-3.840397826549158,-3.840397826549158,-3.840397826549158,-3.840397826549158,-3.840397826549158,-3.840397826549158,-3.840397826549158,-3.840397826549158,-3.840397826549158)
for ($i = 0; $i -lt $r14.Length; $i++) { $ws.Cells.Item(14, 23 + $i).Value = $r14[$i] }
$r15 = @(5.537148498057909,4.858483928622914,4.692141525353577,4.451460339319802,3.379658261193086,-2.491095770678031,-0.236010050592228,-1.162988086281536,0.3986977119751156,0.4839811651348835,0.4839811651348835,0.4839811651348835,0.4839811651348835,0.4839811651348835,0.4839811651348835,0.4839811651348835,0.4839811651348835,0.4839811651348835,0.4839811651348835,0.4839811651348835,0.4839811651348835,0.4839811651348835,0.4839811651348835,0.4839811651348835,0.4839811651348835,0.4839811651348835)
for ($i = 0; $i -lt $r15.Length; $i++) { $ws.Cells.Item(15, 27 + $i).Value = $r15[$i] }
$r16 = @(3.725004516987496,-0.1871946888853837,0.1963989771473607,-0.8351862827690737,1.500453746466346,1.25489130894445,2.177145583294293,0.9587999512773893,1.97557360987699,2.06342951900429,2.06342951900429,2.06342951900429,2.06342951900429,2.06342951900429,2.06342951900429,2.06342951900429,2.06342951900429,2.06342951900429,2.06342951900429,2.06342951900429,2.06342951900429,2.06342951900429)
for ($i = 0; $i -lt $r16.Length; $i++) { $ws.Cells.Item(16, 31 + $i).Value = $r16[$i] }
$r17 = @(-0.6271380860239262,0.4027662762101869,-0.3850568868524862,0.0360848986632023,-1.150174078429844,0.05919149746531627,0.47378583592943,0.3896432785800652,-2.242084520390608,-2.311523918755531,-2.156362896191677,-2.156362896191677,-2.156362896191677,-2.156362896191677,-2.156362896191677,-2.156362896191677,-2.156362896191677,-2.156362896191677,-2.156362896191677)
for ($i = 0; $i -lt $r17.Length; $i++) { $ws.Cells.Item(17, 34 + $i).Value = $r17[$i] }
$r18 = @(-0.8630989168629921,-0.5798629304026526,-0.3296074516580316,4.745563085952331,0.4860776313184267,-0.1671238976421296,0.7251362782769055,-1.437335768580206,-1.379035600217182,-0.8195740704359578,-0.8205034771073372,-0.8205034771073372,-0.8205034771073372,-0.8205034771073372,-0.8205034771073372)
for ($i = 0; $i -lt $r18.Length; $i++) { $ws.Cells.Item(18, 38 + $i).Value = $r18[$i] }
$r19 = @(0.9066518070286467,0.1286052726408693,0.7930386201910888,-0.8900729308547461,-1.053916016632561,-0.4351858173977874,-0.1639718104104904,0.07297157746815053,0.2518138058242414,0.4857399523052974,0.5750555200350504)
for ($i = 0; $i -lt $r19.Length; $i++) { $ws.Cells.Item(19, 42 + $i).Value = $r19[$i] }
$r20 = @(-1.099759890829544,-1.022802765586617,-0.5861676960653317,-1.165535347616076,-1.185914719184011,-0.3387208997876479,0.3466793681147884)
for ($i = 0; $i -lt $r20.Length; $i++) { $ws.Cells.Item(20, 46 + $i).Value = $r20[$i] }
$r21 = @(-1.256374469693167,-0.4360715507411128,0.2338620606656017)
for ($i = 0; $i -lt $r21.Length; $i++) { $ws.Cells.Item(21, 50 + $i).Value = $r21[$i] }
